# Insert 10 blank paragraphs (small 8pt "Garet" font, justified, no
# automatic spacing-after) right after the "Atendiendo a lo establecido..."
# data-protection paragraph and before the existing trailing blank
# paragraph at the end of the document.

$d = $word.ActiveDocument

# Locate the paragraph that starts the well-known clause; this is more
# robust than hard-coding a paragraph index.
$count = $d.Paragraphs.Count
$targetIndex = -1
for ($i = 1; $i -le $count; $i++) {
    $p = $d.Paragraphs.Item($i)
    if ($p.Range.Text -like "Atendiendo a lo establecido*") {
        $targetIndex = $i
    }
}

$target = $d.Paragraphs.Item($targetIndex)
$anchor = $target.Next()

# Insert ten new empty paragraphs immediately before the paragraph that
# currently follows the clause (the pre-existing trailing blank
# paragraph). Each new paragraph inherits that paragraph's formatting
# (Garet/Arial, 8pt, justified, no spacing after), matching the target.
for ($i = 0; $i -lt 10; $i++) {
    $anchor.Range.InsertParagraphBefore()
}
